$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "js copied into" / "jade copied into" columns (K:L) are being removed entirely;
# everything to their right shifts left two columns.
$ws.Range("K1:L1").EntireColumn.Delete()

# Fill in previously-blank "js created" (I) cells for rows where that step is now
# also considered done.
$ws.Range("I8").Value = "x"
$ws.Range("I9").Value = "x"
$ws.Range("I10").Value = "x"
$ws.Range("I13").Value = "x"

# Rows 16 and 17 (Delete Confirmation / Generic Rename) now also have js/jade created.
$ws.Range("I16").Value = "x"
$ws.Range("J16").Value = "x"
$ws.Range("I17").Value = "x"
$ws.Range("J17").Value = "x"

# Update the active selection on the sheet.
$ws.Range("K7").Select()
